# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder: insert "Catar" alphabetically (after Bielorrusia, before
# Corea del Sur). Rather than physically inserting a row, shift the data
# for "Corea del Sur" / "Emiratos Arabes Unidos" down one row and give
# row 36 the new "Catar" figures - matching the target layout exactly.

# Row 36 becomes Catar (new data)
$ws.Range("A36").Value = "Catar"
$ws.Range("B36").Value = 11244
$ws.Range("C36").Value = 957
$ws.Range("D36").Value = 1066
$ws.Range("E36").Value = 10168
$ws.Range("F36").Value = 72
$ws.Range("G36").Value = 0
$ws.Range("H36").Value = 10

# Row 37 becomes Corea del Sur (old row 36 data)
$ws.Range("A37").Value = "Corea del Sur"
$ws.Range("B37").Value = 10738
$ws.Range("C37").Value = 10
$ws.Range("D37").Value = 8764
$ws.Range("E37").Value = 1731
$ws.Range("F37").Value = 55
$ws.Range("G37").Value = 1
$ws.Range("H37").Value = 243

# Row 38 becomes Emiratos Arabes Unidos (old row 37 data)
$ws.Range("A38").Value = "Emiratos Arabes Unidos"
$ws.Range("B38").Value = 10349
$ws.Range("C38").Value = 0
$ws.Range("D38").Value = 1978
$ws.Range("E38").Value = 8295
$ws.Range("F38").Value = 1
$ws.Range("G38").Value = 0
$ws.Range("H38").Value = 76

# --- Plain value corrections (no reordering involved) ---

# Bulgaria (row 84)
$ws.Range("D84").Value = 206
$ws.Range("E84").Value = 1086

# Libano (row 97)
$ws.Range("B97").Value = 710
$ws.Range("C97").Value = 3
$ws.Range("E97").Value = 541

# Malta (row 110)
$ws.Range("B110").Value = 450
$ws.Range("C110").Value = 2
$ws.Range("D110").Value = 286
$ws.Range("E110").Value = 160
$ws.Range("F110").Value = 1

# Kenia (row 117)
$ws.Range("D117").Value = 114
$ws.Range("E117").Value = 235

# Macao (row 169)
$ws.Range("D169").Value = 32
$ws.Range("E169").Value = 13
